$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'269.15"
$ws.Range("G2").Value = "'19"
$ws.Range("D3").Value = "'22.78"
$ws.Range("G3").Value = "'19"
$ws.Range("D4").Value = "'6.340"
$ws.Range("G4").Value = "'19"
$ws.Range("D5").Value = "'0.06175"
$ws.Range("G5").Value = "'19"
$ws.Range("D6").Value = "'3.651"
$ws.Range("G6").Value = "'19"
$ws.Range("D7").Value = "'6.661"
$ws.Range("G7").Value = "'19"
$ws.Range("E8").Value = "7FTXTokenFTTWorstin24h"
$ws.Range("G8").Value = "'19"
$ws.Range("D9").Value = "'0.8281"
$ws.Range("G9").Value = "'19"
$ws.Range("D10").Value = "'0.01371"
$ws.Range("G10").Value = "'19"
$ws.Range("D11").Value = "'0.1604"
$ws.Range("G11").Value = "'19"
$ws.Range("D12").Value = "'0.08287"
$ws.Range("G12").Value = "'19"
$ws.Range("D13").Value = "'0.03556"
$ws.Range("G13").Value = "'19"
$ws.Range("D14").Value = "'0.03199"
$ws.Range("G14").Value = "'19"
$ws.Range("D15").Value = "'0.09324"
$ws.Range("G15").Value = "'19"
$ws.Range("D16").Value = "'3.868"
$ws.Range("G16").Value = "'19"
$ws.Range("G17").Value = "'19"
$ws.Range("D18").Value = "'0.04738"
$ws.Range("G18").Value = "'19"
$ws.Range("D19").Value = "'0.006342"
$ws.Range("G19").Value = "'19"
$ws.Range("G20").Value = "'19"
$ws.Range("D21").Value = "'0.001078"
$ws.Range("G21").Value = "'19"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("G22").Value = "'19"
$ws.Range("D23").Value = "'3.727"
$ws.Range("G23").Value = "'19"
$ws.Range("D24").Value = "'2.412"
$ws.Range("G24").Value = "'19"
$ws.Range("D25").Value = "'0.3304"
$ws.Range("G25").Value = "'19"
$ws.Range("G26").Value = "'19"
$ws.Range("D27").Value = "'0.0002704"
$ws.Range("G27").Value = "'19"
$ws.Range("G28").Value = "'19"
$ws.Range("G29").Value = "'19"
$ws.Range("G30").Value = "'19"
$ws.Range("G31").Value = "'19"
$ws.Range("G32").Value = "'19"
$ws.Range("G33").Value = "'19"
$ws.Range("G34").Value = "'19"
$ws.Range("G35").Value = "'19"
$ws.Range("G36").Value = "'19"
$ws.Range("G37").Value = "'19"
$ws.Range("G38").Value = "'19"
$ws.Range("G39").Value = "'19"
$ws.Range("D40").Value = "'0.04709"
$ws.Range("G40").Value = "'19"
$ws.Range("G41").Value = "'19"
$ws.Range("G42").Value = "'19"
$ws.Range("D43").Value = "'0.003300"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("G43").Value = "'19"
$ws.Range("D44").Value = "'0.01180"
$ws.Range("G44").Value = "'19"
$ws.Range("D45").Value = "'0.00006236"
$ws.Range("G45").Value = "'19"
$ws.Range("D46").Value = "'0.0009899"
$ws.Range("G46").Value = "'19"
$ws.Range("G47").Value = "'19"
$ws.Range("D48").Value = "'0.9200"
$ws.Range("G48").Value = "'19"
$ws.Range("D49").Value = "'0.002351"
$ws.Range("G49").Value = "'19"
$ws.Range("D50").Value = "'0.00001400"
$ws.Range("G50").Value = "'19"
$ws.Range("D51").Value = "'0.01240"
$ws.Range("G51").Value = "'19"
